$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 18:16"

# --- Italia (row 5): updated totals ---
$ws.Range("B5").Value = 69176
$ws.Range("C5").Value = 5249
$ws.Range("D5").Value = 8326
$ws.Range("E5").Value = 54030
$ws.Range("F5").Value = 3393
$ws.Range("G5").Value = 743
$ws.Range("H5").Value = 6820

# --- Estados Unidos (row 6): updated totals ---
$ws.Range("B6").Value = 49594
$ws.Range("C6").Value = 5860
$ws.Range("D6").Value = 361
$ws.Range("E6").Value = 48611
$ws.Range("G6").Value = 69
$ws.Range("H6").Value = 622

# --- Chequia (row 27): updated totals ---
$ws.Range("B27").Value = 1394
$ws.Range("C27").Value = 158
$ws.Range("E27").Value = 1384

# --- Luxemburgo moves up in ranking (was row 33, now row 30) ---
# Row 30 now holds Luxemburgo's new figures
$ws.Range("A30").Value = "Luxemburgo"
$ws.Range("B30").Value = 1099
$ws.Range("C30").Value = 224
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = 1085
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 8

# Ecuador shifts down from row 30 to row 31 (figures unchanged)
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 1049
$ws.Range("C31").Value = 68
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 1019
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 27

# Pakistan shifts down from row 31 to row 32 (figures unchanged)
$ws.Range("A32").Value = "Pakistan"
$ws.Range("B32").Value = 958
$ws.Range("C32").Value = 83
$ws.Range("D32").Value = 13
$ws.Range("E32").Value = 938
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 7

# Chile shifts down from row 32 to row 33 (figures unchanged)
$ws.Range("A33").Value = "Chile"
$ws.Range("B33").Value = 922
$ws.Range("C33").Value = 176
$ws.Range("D33").Value = 17
$ws.Range("E33").Value = 903
$ws.Range("F33").Value = 7
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 2

# --- Bosnia y Herzegovina (row 79): updated totals ---
$ws.Range("B79").Value = 154
$ws.Range("C79").Value = 18
$ws.Range("E79").Value = 150
